$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('赛微电子', '工业富联', '平潭发展'),
    @('实达集团', '赛微电子', '实达集团'),
    @('新 华 都', '航天发展', '航天发展'),
    @('航天发展', '合富中国', '国晟科技'),
    @('特发信息', '特发信息', '孚日股份'),
    @('平潭发展', '东百集团', '特发信息'),
    @('合富中国', '实达集团', '特一药业'),
    @('工业富联', '雷科防务', '合富中国'),
    @('雷科防务', '新 华 都', '新华都'),
    @('中际旭创', '平潭发展', '遥望科技'),
    @('东百集团', '特一药业', '安妮股份'),
    @('孚日股份', '孚日股份', '赛微电子'),
    @('遥望科技', '福日电子', '万通发展'),
    @('特一药业', '茂业商业', '航天动力'),
    @('国晟科技', '永鼎股份', '雷科防务'),
    @('万通发展', '国晟科技', '广济药业'),
    @('航天动力', '安妮股份', '万科A'),
    @('永鼎股份', '航天动力', '国风新材'),
    @('万  科Ａ', '万  科Ａ', '海王生物'),
    @('安妮股份', '中际旭创', '榕基软件')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$wb.Save()